# fix: ajustes na tela de Novo Pedido (checkboxes alinhadas, busca por
# codigo, normalizacao do cilindrico)
#
# The "Pagamentos do Dia" table is reshuffled: the SHINEDUX / BVS 1.61 AR
# pair in rows 3-4 becomes a PEREGO / ULTEX INCOLOR pair, the Base
# 6.00/+2.50 combination is split out into its own row, and the
# BVS 1.67 AR BLUE line is split into two distinct base rows. The table
# grows by two rows and the TOTAL row moves from row 8 down to row 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $fornecedor, $produto, $estoque, $dioptria, $data, $valor) {
    $ws.Range("A$r").Value = $fornecedor
    $ws.Range("B$r").Value = $produto
    $ws.Range("C$r").Value = $estoque
    $ws.Range("D$r").Value = $dioptria
    # keep the date column as plain text instead of letting Excel coerce
    # the "yyyy-mm-dd" string into a real date serial number
    $ws.Range("E$r").Value = "'" + $data
    $ws.Range("F$r").Value = $valor
}

Set-Row 3 "PEREGO"   "ULTEX INCOLOR"    "Não" "Base 4.00 / Adição +2.50" "2025-09-06" 20
Set-Row 4 "PEREGO"   "ULTEX INCOLOR"    "Não" "Base 4.00 / Adição +2.50" "2025-09-06" 20
Set-Row 5 "SHINEDUX" "BVS 1.61 AR"      "Não" "Base 4.00 / Adição +2.50" "2025-09-06" 19
Set-Row 6 "SHINEDUX" "BVS 1.61 AR"      "Não" "Base 6.00 / Adição +2.50" "2025-09-06" 19
Set-Row 7 "SHINEDUX" "BVS 1.67 AR BLUE" "Não" "Base 4.00 / Adição +1.50" "2025-09-06" 80
Set-Row 8 "SHINEDUX" "BVS 1.67 AR BLUE" "Não" "Base 6.00 / Adição +1.50" "2025-09-06" 80

# row 9 stays a blank spacer row

# TOTAL row moves from row 8 to row 10, keeping the bold style
$ws.Range("E10").Value = "TOTAL"
$ws.Range("F10").Value = 243
$ws.Range("E10").Font.Bold = $true
$ws.Range("F10").Font.Bold = $true

# the old row 8 is now an ordinary data row, make sure it is not bold
$ws.Range("A8:F8").Font.Bold = $false

$ws.Range("A1").Select() | Out-Null

Write-Host "Applied Novo Pedido adjustments (rows 3-10)."
